$wb = $excel.ActiveWorkbook

# --- Sheet "Summary" ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.6235955056179775
$ws1.Range("C2").Value = 0.5706638115631691
$ws1.Range("D2").Value = 0.99812734082397
$ws1.Range("E2").Value = 0.726158038147139
$ws1.Range("F2").Value = 0.8680781758957655
$ws1.Range("G2").Value = 0.9701764211705405
$ws1.Range("H2").Value = 0.7825996998134354
$ws1.Range("I2").Value = 533
$ws1.Range("J2").Value = 401
$ws1.Range("K2").Value = 133
$ws1.Range("L2").Value = 1

# --- Sheet "Classification Report" ---
$ws2 = $wb.Worksheets.Item("Classification Report")

# row 2 ("0")
$ws2.Range("B2").Value = 0.9925373134328358
$ws2.Range("C2").Value = 0.249063670411985
$ws2.Range("D2").Value = 0.3982035928143712

# row 3 ("1")
$ws2.Range("B3").Value = 0.5706638115631691
$ws2.Range("C3").Value = 0.99812734082397
$ws2.Range("D3").Value = 0.726158038147139

# row 4 ("accuracy")
$ws2.Range("B4").Value = 0.6235955056179775
$ws2.Range("C4").Value = 0.6235955056179775
$ws2.Range("D4").Value = 0.6235955056179775
$ws2.Range("E4").Value = 0.6235955056179775

# row 5 ("macro avg")
$ws2.Range("B5").Value = 0.7816005624980025
$ws2.Range("C5").Value = 0.6235955056179775
$ws2.Range("D5").Value = 0.562180815480755

# row 6 ("weighted avg")
$ws2.Range("B6").Value = 0.7816005624980026
$ws2.Range("C6").Value = 0.6235955056179775
$ws2.Range("D6").Value = 0.562180815480755

# --- Sheet "Confusion Matrix" ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

# row 2 ("Actual 0")
$ws3.Range("B2").Value = 133
$ws3.Range("C2").Value = 401

# row 3 ("Actual 1")
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 533
